$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Cells.Item(92, 8).Value = 1597.2693
$ws.Cells.Item(92, 9).Value = 663.9524
$ws.Cells.Item(92, 10).Value = 5517.2
$ws.Cells.Item(92, 11).Value = 663.9524
$ws.Cells.Item(92, 12).Value = 5517.2
$ws.Cells.Item(92, 13).Value = 584.0476
$ws.Cells.Item(92, 14).Value = -8013.2

# Row 116
$ws.Cells.Item(116, 8).Value = 6593.5557
$ws.Cells.Item(116, 9).Value = 5766.7144
$ws.Cells.Item(116, 11).Value = 5766.7144
$ws.Cells.Item(116, 13).Value = -2324.7144

# Row 121
$ws.Cells.Item(121, 8).Value = 1500
$ws.Cells.Item(121, 10).Value = 1500
$ws.Cells.Item(121, 12).Value = 4500
$ws.Cells.Item(121, 14).Value = -7994

# Row 132
$ws.Cells.Item(132, 8).Value = 3591770
$ws.Cells.Item(132, 9).Value = 3990606.8
$ws.Cells.Item(132, 11).Value = 11971820.4
$ws.Cells.Item(132, 13).Value = -11969290.4

# Row 137
$ws.Cells.Item(137, 8).Value = 7227.463
$ws.Cells.Item(137, 9).Value = 11527.25
$ws.Cells.Item(137, 11).Value = 34581.75
$ws.Cells.Item(137, 13).Value = -32031.75

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 2433.75
$ws.Cells.Item(2, 9).Value = 2354.5
$ws.Cells.Item(2, 10).Value = 2544.7
$ws.Cells.Item(2, 11).Value = 2354.5
$ws.Cells.Item(2, 12).Value = 2544.7
$ws.Cells.Item(2, 13).Value = -2241.5
$ws.Cells.Item(2, 14).Value = -2770.7

# Row 32
$ws.Cells.Item(32, 8).Value = 17488.426
$ws.Cells.Item(32, 9).Value = 17330.05
$ws.Cells.Item(32, 11).Value = 17330.05
$ws.Cells.Item(32, 13).Value = -17043.05

# Row 33
$ws.Cells.Item(33, 8).Value = 241999.33
$ws.Cells.Item(33, 9).Value = 241999.33
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 241999.33
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = -241670.33
$ws.Cells.Item(33, 14).ClearContents()

# Row 56
$ws.Cells.Item(56, 8).Value = 8500
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 14).ClearContents()

# Row 61
$ws.Cells.Item(61, 8).Value = 3792.6191
$ws.Cells.Item(61, 9).Value = 2926.2354
$ws.Cells.Item(61, 10).Value = 7474.75
$ws.Cells.Item(61, 11).Value = 2926.2354
$ws.Cells.Item(61, 12).Value = 7474.75
$ws.Cells.Item(61, 13).Value = -2714.2354
$ws.Cells.Item(61, 14).Value = -7898.75

# Row 74
$ws.Cells.Item(74, 8).Value = 175971.05
$ws.Cells.Item(74, 9).Value = 180411.36
$ws.Cells.Item(74, 10).Value = 2799
$ws.Cells.Item(74, 11).Value = 180411.36
$ws.Cells.Item(74, 12).Value = 2799
$ws.Cells.Item(74, 13).Value = -179537.36
$ws.Cells.Item(74, 14).Value = -4547

# Row 77
$ws.Cells.Item(77, 8).Value = 175971.05
$ws.Cells.Item(77, 9).Value = 180411.36
$ws.Cells.Item(77, 10).Value = 2799
$ws.Cells.Item(77, 11).Value = 902056.7999999999
$ws.Cells.Item(77, 12).Value = 13995
$ws.Cells.Item(77, 13).Value = -897688.7999999999
$ws.Cells.Item(77, 14).Value = -22731

# Row 116
$ws.Cells.Item(116, 8).Value = 2433.75
$ws.Cells.Item(116, 9).Value = 2354.5
$ws.Cells.Item(116, 10).Value = 2544.7
$ws.Cells.Item(116, 11).Value = 2354.5
$ws.Cells.Item(116, 12).Value = 2544.7
$ws.Cells.Item(116, 13).Value = -60.5
$ws.Cells.Item(116, 14).Value = -7132.7

# Row 136
$ws.Cells.Item(136, 8).Value = 3792.6191
$ws.Cells.Item(136, 9).Value = 2926.2354
$ws.Cells.Item(136, 10).Value = 7474.75
$ws.Cells.Item(136, 11).Value = 8778.706200000001
$ws.Cells.Item(136, 12).Value = 22424.25
$ws.Cells.Item(136, 13).Value = -6228.706200000001
$ws.Cells.Item(136, 14).Value = -27524.25

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 2433.75
$ws.Cells.Item(3, 9).Value = 2354.5
$ws.Cells.Item(3, 10).Value = 2544.7
$ws.Cells.Item(3, 11).Value = 2354.5
$ws.Cells.Item(3, 12).Value = 2544.7
$ws.Cells.Item(3, 13).Value = -2240.5
$ws.Cells.Item(3, 14).Value = -2772.7

# Row 75
$ws.Cells.Item(75, 8).Value = 28747.125
$ws.Cells.Item(75, 10).Value = 34996.332
$ws.Cells.Item(75, 12).Value = 34996.332
$ws.Cells.Item(75, 14).Value = -36868.332

# Row 78
$ws.Cells.Item(78, 8).Value = 28747.125
$ws.Cells.Item(78, 10).Value = 34996.332
$ws.Cells.Item(78, 12).Value = 104988.996
$ws.Cells.Item(78, 14).Value = -114348.996

# Row 94
$ws.Cells.Item(94, 8).Value = 1021.907
$ws.Cells.Item(94, 9).Value = 1099.3125
$ws.Cells.Item(94, 11).Value = 1099.3125
$ws.Cells.Item(94, 13).Value = -648.3125

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 4352524.5
$ws.Cells.Item(31, 9).Value = 5559433.5
$ws.Cells.Item(31, 10).Value = 7652.4
$ws.Cells.Item(31, 11).Value = 5559433.5
$ws.Cells.Item(31, 12).Value = 7652.4
$ws.Cells.Item(31, 13).Value = -5559138.5
$ws.Cells.Item(31, 14).Value = -8242.4

# Row 34
$ws.Cells.Item(34, 8).Value = 4352524.5
$ws.Cells.Item(34, 9).Value = 5559433.5
$ws.Cells.Item(34, 10).Value = 7652.4
$ws.Cells.Item(34, 11).Value = 5559433.5
$ws.Cells.Item(34, 12).Value = 7652.4
$ws.Cells.Item(34, 13).Value = -5559231.5
$ws.Cells.Item(34, 14).Value = -8056.4

# Row 62
$ws.Cells.Item(62, 8).Value = 10094.5
$ws.Cells.Item(62, 9).Value = 4711.25
$ws.Cells.Item(62, 11).Value = 4711.25
$ws.Cells.Item(62, 13).Value = -4087.25

# Row 65
$ws.Cells.Item(65, 8).Value = 10094.5
$ws.Cells.Item(65, 9).Value = 4711.25
$ws.Cells.Item(65, 11).Value = 23556.25
$ws.Cells.Item(65, 13).Value = -20436.25

# Row 94
$ws.Cells.Item(94, 8).Value = 1401.875
$ws.Cells.Item(94, 10).Value = 1105.4
$ws.Cells.Item(94, 12).Value = 1105.4
$ws.Cells.Item(94, 14).Value = -2007.4

# Row 122
$ws.Cells.Item(122, 8).Value = 68150
$ws.Cells.Item(122, 9).Value = 68150
$ws.Cells.Item(122, 11).Value = 204450
$ws.Cells.Item(122, 13).Value = -202000

# Row 132
$ws.Cells.Item(132, 8).Value = 39207.5
$ws.Cells.Item(132, 9).Value = 44268.93
$ws.Cells.Item(132, 11).Value = 132806.79
$ws.Cells.Item(132, 13).Value = -130276.79

# Row 141
$ws.Cells.Item(141, 8).Value = 200862.44
$ws.Cells.Item(141, 10).Value = 200862.44
$ws.Cells.Item(141, 12).Value = 200862.44
$ws.Cells.Item(141, 14).Value = -211222.44

$ws = $wb.Worksheets.Item("GSM")
# Row 49
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).ClearContents()

# Row 55
$ws.Cells.Item(55, 8).Value = 10628.875
$ws.Cells.Item(55, 9).Value = 7999.5
$ws.Cells.Item(55, 10).Value = 11505.333
$ws.Cells.Item(55, 11).Value = 7999.5
$ws.Cells.Item(55, 12).Value = 11505.333
$ws.Cells.Item(55, 13).Value = -7672.5
$ws.Cells.Item(55, 14).Value = -12159.333

# Row 70
$ws.Cells.Item(70, 8).Value = 6371.136
$ws.Cells.Item(70, 9).Value = 6341.778
$ws.Cells.Item(70, 10).Value = 6391.4614
$ws.Cells.Item(70, 11).Value = 6341.778
$ws.Cells.Item(70, 12).Value = 6391.4614
$ws.Cells.Item(70, 13).Value = -6071.778
$ws.Cells.Item(70, 14).Value = -6931.4614

# Row 73
$ws.Cells.Item(73, 8).Value = 6371.136
$ws.Cells.Item(73, 9).Value = 6341.778
$ws.Cells.Item(73, 10).Value = 6391.4614
$ws.Cells.Item(73, 11).Value = 6341.778
$ws.Cells.Item(73, 12).Value = 6391.4614
$ws.Cells.Item(73, 13).Value = -5405.778
$ws.Cells.Item(73, 14).Value = -8263.4614

# Row 80
$ws.Cells.Item(80, 8).Value = 4796.0454
$ws.Cells.Item(80, 9).Value = 2969.6875
$ws.Cells.Item(80, 11).Value = 2969.6875
$ws.Cells.Item(80, 13).Value = -1971.6875

# Row 83
$ws.Cells.Item(83, 8).Value = 4796.0454
$ws.Cells.Item(83, 9).Value = 2969.6875
$ws.Cells.Item(83, 11).Value = 14848.4375
$ws.Cells.Item(83, 13).Value = -9856.4375

# Row 102
$ws.Cells.Item(102, 8).Value = 19384.768
$ws.Cells.Item(102, 9).Value = 22981.584
$ws.Cells.Item(102, 10).Value = 4997.5
$ws.Cells.Item(102, 11).Value = 22981.584
$ws.Cells.Item(102, 12).Value = 4997.5
$ws.Cells.Item(102, 13).Value = -21359.584
$ws.Cells.Item(102, 14).Value = -8241.5

# Row 113
$ws.Cells.Item(113, 8).Value = 1969.1177
$ws.Cells.Item(113, 9).Value = 1966.0667
$ws.Cells.Item(113, 11).Value = 1966.0667
$ws.Cells.Item(113, 13).Value = 203.9332999999999

# Row 122
$ws.Cells.Item(122, 8).Value = 4012.2778
$ws.Cells.Item(122, 9).Value = 3732.875
$ws.Cells.Item(122, 10).Value = 6247.5
$ws.Cells.Item(122, 11).Value = 11198.625
$ws.Cells.Item(122, 12).Value = 18742.5
$ws.Cells.Item(122, 13).Value = -8748.625
$ws.Cells.Item(122, 14).Value = -23642.5

# Row 126
$ws.Cells.Item(126, 8).Value = 2594.5
$ws.Cells.Item(126, 9).Value = 2482.8635
$ws.Cells.Item(126, 11).Value = 7448.5905
$ws.Cells.Item(126, 13).Value = -4978.5905

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 3467.75
$ws.Cells.Item(7, 9).Value = 3549.4
$ws.Cells.Item(7, 10).Value = 3331.6667
$ws.Cells.Item(7, 11).Value = 3549.4
$ws.Cells.Item(7, 12).Value = 3331.6667
$ws.Cells.Item(7, 13).Value = -3437.4
$ws.Cells.Item(7, 14).Value = -3555.6667

# Row 42
$ws.Cells.Item(42, 8).Value = 17999.5
$ws.Cells.Item(42, 10).Value = 17999.5
$ws.Cells.Item(42, 12).Value = 17999.5
$ws.Cells.Item(42, 14).Value = -19125.5

# Row 49
$ws.Cells.Item(49, 8).Value = 17999.5
$ws.Cells.Item(49, 10).Value = 17999.5
$ws.Cells.Item(49, 12).Value = 17999.5
$ws.Cells.Item(49, 14).Value = -18293.5

# Row 55
$ws.Cells.Item(55, 8).Value = 1899.1786
$ws.Cells.Item(55, 10).Value = 3043
$ws.Cells.Item(55, 12).Value = 3043
$ws.Cells.Item(55, 14).Value = -3389

# Row 126
$ws.Cells.Item(126, 8).Value = 3467.75
$ws.Cells.Item(126, 9).Value = 3549.4
$ws.Cells.Item(126, 10).Value = 3331.6667
$ws.Cells.Item(126, 11).Value = 10648.2
$ws.Cells.Item(126, 12).Value = 9995.000100000001
$ws.Cells.Item(126, 13).Value = -8178.200000000001
$ws.Cells.Item(126, 14).Value = -14935.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Cells.Item(113, 8).Value = 853.5
$ws.Cells.Item(113, 9).Value = 521.38464
$ws.Cells.Item(113, 10).Value = 1333.2222
$ws.Cells.Item(113, 11).Value = 1564.15392
$ws.Cells.Item(113, 12).Value = 3999.6666
$ws.Cells.Item(113, 13).Value = 605.84608
$ws.Cells.Item(113, 14).Value = -8339.6666

# Row 141
$ws.Cells.Item(141, 8).Value = 98735.625
$ws.Cells.Item(141, 10).Value = 98735.625
$ws.Cells.Item(141, 12).Value = 98735.625
$ws.Cells.Item(141, 14).Value = -109095.625
